$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking strings (e.g. "235.49") that must
# stay plain text, matching the workbook's existing inline-string cells, instead of
# being auto-converted to numbers. Temporarily force text format on the column,
# write the values, then restore the default "Normal" style so no stray per-cell
# style index is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Column D (Price) updates
$ws.Range("D2").Value = '24.954.55'
$ws.Range("D3").Value = '1.640.47'
$ws.Range("D4").Value = '0.9988'
$ws.Range("D5").Value = '235.49'
$ws.Range("D7").Value = '0.4778'
$ws.Range("D8").Value = '0.2571'
$ws.Range("D9").Value = '0.05991'
$ws.Range("D10").Value = '0.07213'
$ws.Range("D11").Value = '1.643.90'
$ws.Range("D13").Value = '0.6165'
$ws.Range("D14").Value = '4.486'
$ws.Range("D15").Value = '72.65'
$ws.Range("D16").Value = '1.001'
$ws.Range("D17").Value = '0.9989'
$ws.Range("D18").Value = '24.943.05'
$ws.Range("D19").Value = '11.29'
$ws.Range("D20").Value = '0.000006592'
$ws.Range("D21").Value = '4.454'
$ws.Range("D22").Value = '1.852.51'
$ws.Range("D23").Value = '8.576'
$ws.Range("D24").Value = '5.274'
$ws.Range("D25").Value = '132.08'
$ws.Range("D26").Value = '14.82'
$ws.Range("D27").Value = '1.382'
$ws.Range("D28").Value = '103.03'
$ws.Range("D29").Value = '1.659'
$ws.Range("D31").Value = '0.07797'
$ws.Range("D32").Value = '3.535'
$ws.Range("D33").Value = '0.04409'
$ws.Range("D34").Value = '0.9998'
$ws.Range("D35").Value = '2.591'
$ws.Range("D36").Value = '0.9247'
$ws.Range("D37").Value = '0.5831'
$ws.Range("D38").Value = '2.554'
$ws.Range("D39").Value = '0.01561'
$ws.Range("D40").Value = '0.8408'
$ws.Range("D41").Value = '0.9991'
$ws.Range("D42").Value = '1.798'
$ws.Range("D43").Value = '97.21'
$ws.Range("D44").Value = '0.3710'
$ws.Range("D45").Value = '4.752'
$ws.Range("D46").Value = '0.1149'
$ws.Range("D47").Value = '6.079'
$ws.Range("D48").Value = '0.05194'
$ws.Range("D49").Value = '29.63'
$ws.Range("D50").Value = '0.9993'

$priceRange.Style = "Normal"

# Column E (Volume(1h)) updates
$ws.Range("E2").Value = '  -3.67%  '
$ws.Range("E3").Value = '  -5.78%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("E5").Value = '  -5.09%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("E7").Value = '  -6.55%  '
$ws.Range("E8").Value = '  -6.19%  '
$ws.Range("E9").Value = '  -3.12%  '
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("E11").Value = '  -5.54%  '
$ws.Range("E12").Value = '  -2.19%  '
$ws.Range("E13").Value = '  -4.80%  '
$ws.Range("E14").Value = '  -3.09%  '
$ws.Range("E15").Value = '  -6.40%  '
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("E18").Value = '  -3.80%  '
$ws.Range("E19").Value = '  -4.36%  '
$ws.Range("E20").Value = '  -3.17%  '
$ws.Range("E21").Value = '  +4.21%  '
$ws.Range("E22").Value = '  -5.63%  '
$ws.Range("E23").Value = '  -0.78%  '
$ws.Range("E24").Value = '  -2.10%  '
$ws.Range("E25").Value = '  -3.19%  '
$ws.Range("E26").Value = '  -2.68%  '
$ws.Range("E27").Value = '  -8.14%  '
$ws.Range("E28").Value = '  -2.22%  '
$ws.Range("E29").Value = '  -6.37%  '
$ws.Range("E30").Value = '  -4.74%  '
$ws.Range("E31").Value = '  -5.15%  '
$ws.Range("E32").Value = '  -2.92%  '
$ws.Range("E33").Value = '  -5.76%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("E35").Value = '  -2.51%  '
$ws.Range("E36").Value = '  -7.36%  '
$ws.Range("E37").Value = '  -6.87%  '
$ws.Range("E38").Value = '  -6.39%  '
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("E40").Value = '  +10.90%  '
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("E42").Value = '  -6.44%  '
$ws.Range("E43").Value = '  -2.79%  '
$ws.Range("E44").Value = '  -3.55%  '
$ws.Range("E45").Value = '  -4.75%  '
$ws.Range("E46").Value = '  +1.64%  '
$ws.Range("E47").Value = '  -3.34%  '
$ws.Range("E48").Value = '  -0.79%  '
$ws.Range("E49").Value = '  -3.29%  '
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("E51").Value = '  -0.54%  '
